$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B337").Value = 2019
$ws.Range("B359").Value = 2019
$ws.Range("B360").Value = 2019
$ws.Range("B403").Value = 2019
$ws.Range("B442").Value = 2019
$ws.Range("B1161").Value = 44834
$ws.Range("B1162").Value = 44832
$ws.Range("B1164").Value = 44834
$ws.Range("B1165").Value = 44822
$ws.Range("B1166").Value = 44815
$ws.Range("B1168").Value = 44836
$ws.Range("B1170").Value = 44820
$ws.Range("B1171").Value = 44833
$ws.Range("B1172").Value = 44835
$ws.Range("B1174").Value = 44822
$ws.Range("B1175").Value = 44831
$ws.Range("B1176").Value = 44822
$ws.Range("B1179").Value = 44835
$ws.Range("B1180").Value = 44822
$ws.Range("B1181").Value = 44836
$ws.Range("B1182").Value = 44834
$ws.Range("B1184").Value = 44822
$ws.Range("B1185").Value = 44834
$ws.Range("B1188").Value = 44832
$ws.Range("B1189").Value = 44830
$ws.Range("B1190").Value = 44813
$ws.Range("B1191").Value = 44829
$ws.Range("B1192").Value = 44815
$ws.Range("B1193").Value = 44822
$ws.Range("B1194").Value = 44834
$ws.Range("B1195").Value = 44830
$ws.Range("B1196").Value = 44834
$ws.Range("B1198").Value = 44822
$ws.Range("B1199").Value = 44815
$ws.Range("B1200").Value = 44822
$ws.Range("B1205").Value = 44822
$ws.Range("B1206").Value = 44820
$ws.Range("B1207").Value = 44834
$ws.Range("B1208").Value = 44834
$ws.Range("B1210").Value = 44831
$ws.Range("B1211").Value = 44836
$ws.Range("B1212").Value = 44832
$ws.Range("B1213").Value = 44831
$ws.Range("B1214").Value = 44834
$ws.Range("B1215").Value = 44833
$ws.Range("B1216").Value = 44835
$ws.Range("B1219").Value = 44829
$ws.Range("B1220").Value = 44832
$ws.Range("B1221").Value = 44830
$ws.Range("B1223").Value = 44832
$ws.Range("B1224").Value = 44830
$ws.Range("B1226").Value = 44833
$ws.Range("B1228").Value = 44822
$ws.Range("B1232").Value = 44822
$ws.Range("B1234").Value = 44822
$ws.Range("B1236").Value = 44822
$ws.Range("B1238").Value = 44836
$ws.Range("B1239").Value = 44834
$ws.Range("B1241").Value = 44834
$ws.Range("B1243").Value = 44836
$ws.Range("B1244").Value = 44834
$ws.Range("B1245").Value = 44835
$ws.Range("B1246").Value = 44820
$ws.Range("B1247").Value = 44831
$ws.Range("B1248").Value = 44830
$ws.Range("B1250").Value = 44836
$ws.Range("B1251").Value = 44832
$ws.Range("B1252").Value = 44830
$ws.Range("B1253").Value = 44826
$ws.Range("B1255").Value = 44836
$ws.Range("B1256").Value = 44836
$ws.Range("B1257").Value = 44834
$ws.Range("B1258").Value = 44829
$ws.Range("B1260").Value = 44836
$ws.Range("B1261").Value = 44834
$ws.Range("B1262").Value = 44822
$ws.Range("B1263").Value = 44836
$ws.Range("B1264").Value = 44820
$ws.Range("B1267").Value = 44833
$ws.Range("B1268").Value = 44832
$ws.Range("B1270").Value = 44830
$ws.Range("B1271").Value = 44822
$ws.Range("B1272").Value = 44825
$ws.Range("B1274").Value = 44827
$ws.Range("B1275").Value = 44829
$ws.Range("B1277").Value = 44836
$ws.Range("B1278").Value = 44829
$ws.Range("B1280").Value = 44819
$ws.Range("B1281").Value = 44832
$ws.Range("B1283").Value = 44822
$ws.Range("B1284").Value = 44822
$ws.Range("B1285").Value = 44830
$ws.Range("B1286").Value = 44820
$ws.Range("B1287").Value = 44815
$ws.Range("B1288").Value = 44822
$ws.Range("B1289").Value = 44835
$ws.Range("B1293").Value = 44822
$ws.Range("B1295").Value = 44820
$ws.Range("B1296").Value = 44822
$ws.Range("B1297").Value = 44822
$ws.Range("B1298").Value = 44836
$ws.Range("B1299").Value = 44822
$ws.Range("B1300").Value = 44830
$ws.Range("B1301").Value = 44815
$ws.Range("B1302").Value = 44822
$ws.Range("B1303").Value = 44834
$ws.Range("B1306").Value = 44833
$ws.Range("B1307").Value = 44829
$ws.Range("B1309").Value = 44831
$ws.Range("B1310").Value = 44831
$ws.Range("B1311").Value = 44834
$ws.Range("B1312").Value = 44834
$ws.Range("B1314").Value = 44835
$ws.Range("B1315").Value = 44819
$ws.Range("B1316").Value = 44830
$ws.Range("B1317").Value = 44835
$ws.Range("B1318").Value = 44820
$ws.Range("B1319").Value = 44834
$ws.Range("B1322").Value = 44832
$ws.Range("B1323").Value = 44822
$ws.Range("B1324").Value = 44834
$ws.Range("B1325").Value = 44822
$ws.Range("B1326").Value = 44831
$ws.Range("B1328").Value = 44822
$ws.Range("B1329").Value = 44832
$ws.Range("B1331").Value = 44830
$ws.Range("B1332").Value = 44822
$ws.Range("B1335").Value = 44831
$ws.Range("B1337").Value = 44822
$ws.Range("B1338").Value = 44822
$ws.Range("B1342").Value = 44833
$ws.Range("B1343").Value = 44815
$ws.Range("B1344").Value = 44820
$ws.Range("B1345").Value = 44822
$ws.Range("B1346").Value = 44826
$ws.Range("B1349").Value = 44815
$ws.Range("B1350").Value = 44834
$ws.Range("B1354").Value = 44832
$ws.Range("B1356").Value = 44825
$ws.Range("B1358").Value = 44836
$ws.Range("B1360").Value = 44836
$ws.Range("B1361").Value = 44822
$ws.Range("B1362").Value = 44822
$ws.Range("B1364").Value = 44834
$ws.Range("B1365").Value = 44831
$ws.Range("B1367").Value = 44820
$ws.Range("B1369").Value = 44820
$ws.Range("B1370").Value = 44819
$ws.Range("B1373").Value = 44830
$ws.Range("B1374").Value = 44832
$ws.Range("B1375").Value = 44834
